# Auto-generated edit script: apply numeric corrections to Leve profit sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across 8 job sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("L21").Value = 40000
$ws.Range("N21").Value = -40936
$ws.Range("M21").Value = -47876.668
$ws.Range("K21").Value = 48344.668
$ws.Range("I21").Value = 48344.668
$ws.Range("H21").Value = 45006.8
$ws.Range("J21").Value = 40000
$ws.Range("L23").Value = 40000
$ws.Range("N23").Value = -40468
$ws.Range("M23").Value = -48110.668
$ws.Range("K23").Value = 48344.668
$ws.Range("I23").Value = 48344.668
$ws.Range("H23").Value = 45006.8
$ws.Range("J23").Value = 40000
$ws.Range("M34").Value = -1294.4445
$ws.Range("K34").Value = 1497.4445
$ws.Range("I34").Value = 1497.4445
$ws.Range("H34").Value = 1497.4445
$ws.Range("M36").Value = -782.4445000000001
$ws.Range("K36").Value = 1497.4445
$ws.Range("I36").Value = 1497.4445
$ws.Range("H36").Value = 1497.4445
$ws.Range("L88").Value = 3759.7778
$ws.Range("N88").Value = -4571.7778
$ws.Range("M88").Value = -369.5
$ws.Range("K88").Value = 775.5
$ws.Range("I88").Value = 775.5
$ws.Range("H88").Value = 2841.5386
$ws.Range("J88").Value = 3759.7778
$ws.Range("L91").Value = 3759.7778
$ws.Range("N91").Value = -6567.7778
$ws.Range("M91").Value = 628.5
$ws.Range("K91").Value = 775.5
$ws.Range("I91").Value = 775.5
$ws.Range("H91").Value = 2841.5386
$ws.Range("J91").Value = 3759.7778
$ws.Range("L105").Value = 64000
$ws.Range("H105").Value = 64000
$ws.Range("J105").Value = 64000
$ws.Range("N105").Value = -70988
$ws.Range("L132").Value = 34372.66800000001
$ws.Range("N132").Value = -39432.66800000001
$ws.Range("M132").Value = -994320.5
$ws.Range("K132").Value = 996850.5
$ws.Range("I132").Value = 332283.5
$ws.Range("H132").Value = 276755.94
$ws.Range("J132").Value = 11457.556
$ws.Range("M141").Value = -782.1538
$ws.Range("K141").Value = 5962.1538
$ws.Range("I141").Value = 1987.3846
$ws.Range("H141").Value = 2739.25

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M132").Value = -6590027
$ws.Range("K132").Value = 6592557
$ws.Range("I132").Value = 2197519
$ws.Range("H132").Value = 1587444.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M99").Value = -5682.8477
$ws.Range("K99").Value = 7180.8477
$ws.Range("I99").Value = 7180.8477
$ws.Range("H99").Value = 7322.79
$ws.Range("L140").Value = 78999.91
$ws.Range("N140").Value = -89359.91
$ws.Range("H140").Value = 78999.91
$ws.Range("J140").Value = 78999.91

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("L16").Value = 12496.833
$ws.Range("N16").Value = -13070.833
$ws.Range("M16").Value = -16668066
$ws.Range("K16").Value = 16668353
$ws.Range("I16").Value = 16668353
$ws.Range("H16").Value = 13892377
$ws.Range("J16").Value = 12496.833
$ws.Range("M31").Value = -13925.467
$ws.Range("K31").Value = 14220.467
$ws.Range("I31").Value = 14220.467
$ws.Range("H31").Value = 9234.212
$ws.Range("M34").Value = -14018.467
$ws.Range("K34").Value = 14220.467
$ws.Range("I34").Value = 14220.467
$ws.Range("H34").Value = 9234.212
$ws.Range("L57").Value = 40000
$ws.Range("N57").Value = -41120
$ws.Range("H57").Value = 40000
$ws.Range("J57").Value = 40000
$ws.Range("L113").Value = 12496.833
$ws.Range("N113").Value = -16836.833
$ws.Range("M113").Value = -16666183
$ws.Range("K113").Value = 16668353
$ws.Range("I113").Value = 16668353
$ws.Range("H113").Value = 13892377
$ws.Range("J113").Value = 12496.833
$ws.Range("M134").Value = -15944.1
$ws.Range("K134").Value = 18479.1
$ws.Range("I134").Value = 6159.7
$ws.Range("H134").Value = 11160.789

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("L25").Value = 0
$ws.Range("K25").Value = 6003
$ws.Range("I25").Value = 2001
$ws.Range("H25").Value = 2001
$ws.Range("J25").Value = 0
$ws.Range("M25").Value = -5834
$ws.Range("M29").Value = -6986.8568
$ws.Range("K29").Value = 7263.8568
$ws.Range("I29").Value = 2421.2856
$ws.Range("H29").Value = 1938.6666
$ws.Range("L30").Value = 0
$ws.Range("K30").Value = 6003
$ws.Range("I30").Value = 2001
$ws.Range("H30").Value = 2001
$ws.Range("J30").Value = 0
$ws.Range("M30").Value = -5901
$ws.Range("L37").Value = 342564
$ws.Range("N37").Value = -342788
$ws.Range("H37").Value = 114188
$ws.Range("J37").Value = 114188
$ws.Range("N25").ClearContents()
$ws.Range("N30").ClearContents()

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("L80").Value = 12074.9
$ws.Range("N80").Value = -14070.9
$ws.Range("M80").Value = -4389
$ws.Range("K80").Value = 5387
$ws.Range("I80").Value = 5387
$ws.Range("H80").Value = 10164.071
$ws.Range("J80").Value = 12074.9
$ws.Range("L83").Value = 60374.5
$ws.Range("N83").Value = -70358.5
$ws.Range("M83").Value = -21943
$ws.Range("K83").Value = 26935
$ws.Range("I83").Value = 5387
$ws.Range("H83").Value = 10164.071
$ws.Range("J83").Value = 12074.9
$ws.Range("L132").Value = 20016.6
$ws.Range("N132").Value = -25076.6
$ws.Range("M132").Value = -14435.2139
$ws.Range("K132").Value = 16965.2139
$ws.Range("I132").Value = 5655.0713
$ws.Range("H132").Value = 6009.884
$ws.Range("J132").Value = 6672.2

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("M55").Value = -3842.182
$ws.Range("K55").Value = 4015.182
$ws.Range("I55").Value = 4015.182
$ws.Range("H55").Value = 4861.3237
$ws.Range("L68").Value = 949.5
$ws.Range("N68").Value = -2447.5
$ws.Range("M68").Value = -1190.3846
$ws.Range("K68").Value = 1939.3846
$ws.Range("I68").Value = 1939.3846
$ws.Range("H68").Value = 1807.4
$ws.Range("J68").Value = 949.5
$ws.Range("L71").Value = 4747.5
$ws.Range("N71").Value = -12235.5
$ws.Range("M71").Value = -5952.923000000001
$ws.Range("K71").Value = 9696.923000000001
$ws.Range("I71").Value = 1939.3846
$ws.Range("H71").Value = 1807.4
$ws.Range("J71").Value = 949.5
$ws.Range("L106").Value = 40000
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("N106").Value = -42524
$ws.Range("M132").Value = -9039.3572
$ws.Range("K132").Value = 11569.3572
$ws.Range("I132").Value = 3856.4524
$ws.Range("H132").Value = 5082.373

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K26").Value = 12506
$ws.Range("I26").Value = 12506
$ws.Range("H26").Value = 12506
$ws.Range("M26").Value = -12213
$ws.Range("L96").Value = 3110
$ws.Range("N96").Value = -5856
$ws.Range("M96").Value = -1276
$ws.Range("K96").Value = 2649
$ws.Range("I96").Value = 2649
$ws.Range("H96").Value = 2912.4285
$ws.Range("J96").Value = 3110
$ws.Range("L135").Value = 88178.39999999999
$ws.Range("N135").Value = -98318.39999999999
$ws.Range("H135").Value = 88178.39999999999
$ws.Range("J135").Value = 88178.39999999999
